$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being stored as text (not auto-converted to
# numbers by Excel) for every data row, matching the original inlineStr cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "76.333.96"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.040.81"
$ws.Range("E3").Value = "  +3.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - Solana
$ws.Range("D5").Value = "198.88"
$ws.Range("E5").Value = "  -2.87%  "

# Row 6 - BNB
$ws.Range("D6").Value = "620.80"
$ws.Range("E6").Value = "  +3.41%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.548"
$ws.Range("E8").Value = "  -1.11%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.205"
$ws.Range("E9").Value = "  +3.02%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.039.61"
$ws.Range("E10").Value = "  +3.16%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.438"
$ws.Range("E11").Value = "  -1.23%  "

# Row 13 - Toncoin
$ws.Range("D13").Value = "5.25"
$ws.Range("E13").Value = "  +6.20%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.598.33"
$ws.Range("E14").Value = "  +3.04%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "28.86"
$ws.Range("E15").Value = "  +1.93%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "76.301.65"
$ws.Range("E16").Value = "  +0.25%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.0000193"
$ws.Range("E17").Value = "  +0.86%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.035.06"
$ws.Range("E18").Value = "  +2.66%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "13.46"
$ws.Range("E19").Value = "  +1.17%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "8.95"
$ws.Range("E20").Value = "  +0.95%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "378.82"
$ws.Range("E21").Value = "  +1.21%  "

# Row 22 - SuiNetwork
$ws.Range("D22").Value = "2.30"
$ws.Range("E22").Value = "  -1.04%  "

# Row 23 - Polkadot
$ws.Range("D23").Value = "4.35"
$ws.Range("E23").Value = "  +0.52%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "73.09"
$ws.Range("E24").Value = "  +1.68%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.11%  "

# Row 27 - NEARProtocol
$ws.Range("D27").Value = "4.34"
$ws.Range("E27").Value = "  -0.33%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "9.68"
$ws.Range("E28").Value = "  -0.44%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0000109"
$ws.Range("E29").Value = "  -0.72%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.19%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "8.26"
$ws.Range("E31").Value = "  +4.52%  "

# Row 32 - Fetch.AI
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -0.43%  "

# Row 33 - PancakeSwap
$ws.Range("D33").Value = "1.94"
$ws.Range("E33").Value = "  +4.21%  "

# Row 34 - Bittensor
$ws.Range("D34").Value = "490.54"
$ws.Range("E34").Value = "  -3.18%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "20.59"
$ws.Range("E36").Value = "  +1.15%  "

# Row 37 - Monero
$ws.Range("D37").Value = "162.74"
$ws.Range("E37").Value = "  -0.57%  "

# Row 38 - WhiteBITCoin
$ws.Range("D38").Value = "20.03"
$ws.Range("E38").Value = "  +1.90%  "

# Row 39 - now Kaspa (was PolygonEcosystemToken)
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.116"
$ws.Range("E39").Value = "  +2.80%  "

# Row 40 - now PolygonEcosystemToken (was Kaspa)
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "0.381"
$ws.Range("E40").Value = "  +1.65%  "

# Row 41 - Aave
$ws.Range("D41").Value = "189.46"
$ws.Range("E41").Value = "  +4.18%  "

# Row 42 - Cronos
$ws.Range("D42").Value = "0.104"
$ws.Range("E42").Value = "  -5.30%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  -0.05%  "

# Row 44 - Mantle
$ws.Range("D44").Value = "0.803"
$ws.Range("E44").Value = "  +20.64%  "

# Row 45 - RenderToken
$ws.Range("D45").Value = "5.08"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46 - ImmutableX
$ws.Range("D46").Value = "1.26"
$ws.Range("E46").Value = "  +4.22%  "

# Row 47 - OKB
$ws.Range("D47").Value = "41.73"
$ws.Range("E47").Value = "  +3.78%  "

# Row 48 - Stacks
$ws.Range("D48").Value = "1.64"
$ws.Range("E48").Value = "  -1.87%  "

# Row 49 - dogwifhat
$ws.Range("D49").Value = "2.41"
$ws.Range("E49").Value = "  +1.89%  "

# Row 50 - ARBITRUM
$ws.Range("D50").Value = "0.602"
$ws.Range("E50").Value = "  +3.03%  "

# Row 51 - Filecoin
$ws.Range("D51").Value = "3.88"
$ws.Range("E51").Value = "  +2.46%  "
